$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.323.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.92%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.828.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.59%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -1.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'314.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -1.11%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4279"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3700"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.68%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.13%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8657"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.27%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'21.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.10%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.814.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.58%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.721"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.87%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.07109"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.01%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -3.37%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'89.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.26%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.25%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008876"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'15.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.88%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'27.331.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.151"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.60%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.17%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.039.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.53%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -1.24%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'153.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.43%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.54%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.157"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.33%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.258"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'116.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.17%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08913"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.75%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.56%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7608"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.20%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.465"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.69%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.841"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.17%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.113"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.13%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01983"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.49%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05291"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.74%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'7.207"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.54%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.61%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1699"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.76%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.5059"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.89%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.706"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.21%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'10.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.68%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'107.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.76%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.4769"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.13%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.14%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'NEARProtocol"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.671"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.66%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.06374"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.34%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.846"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.90%  "
$ws.Range("E51").Style = "Normal"
